# Add 2022-Q1 data.
#
# Before:  sheets = [ 2021-Q4 (detail), 总计 (totals) ]
# After:   sheets = [ 2021-Q4 (detail), 2022-Q1 (detail, NEW), 总计 (totals, now with 2 rows) ]
#
# The existing "总计" worksheet is repurposed into the new "2022-Q1" detail
# sheet (same layout/format as "2021-Q4"), and a brand new "总计" worksheet
# is appended right after it, holding the (now two-row) totals table.

$wb = $excel.ActiveWorkbook

$detail2021 = $wb.Worksheets.Item("2021-Q4")
$totals = $wb.Worksheets.Item("总计")

# A guaranteed-blank, default-styled cell (outside any table) used below to
# strip stray number-formatting off freshly-written text cells.
$blank = $detail2021.Range("A1")

# Snapshot the current totals rows (header + 2021-Q4 row) before the "总计"
# sheet gets repurposed, so the new totals sheet can be rebuilt from them.
# NOTE: in this host, the ".Value" GETTER is unreliable (it returns a
# reflection descriptor instead of the cell's value) - ".Value2" works.
# ".Value" as a SETTER works fine, so it is still used for writes below.
$totalsHeaderB = $totals.Range("B1").Value2
$totalsHeaderC = $totals.Range("C1").Value2
$totalsHeaderD = $totals.Range("D1").Value2
$oldTotalsDate = $totals.Range("B2").Value2
$oldTotalsCount = $totals.Range("C2").Value2
$oldTotalsValue = $totals.Range("D2").Value2

# --- Turn the old "总计" sheet into the new "2022-Q1" detail sheet ---------

# Copy the full formatting (styles/borders/fonts/alignment) of the detail
# table from "2021-Q4" onto "总计" so the new sheet matches the established
# layout, without touching "总计"'s own values yet.
$detail2021.Range("A1:H2").Copy()
$totals.Range("A1:H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totals.Name = "2022-Q1"
$newDetail = $totals

$newDetail.Range("B1").Value = "基金代码"
$newDetail.Range("C1").Value = "基金名称"
$newDetail.Range("D1").Value = "基金规模"
$newDetail.Range("E1").Value = "股票总仓位"
$newDetail.Range("F1").Value = "仓位占比"
$newDetail.Range("G1").Value = "持有市值(亿元)"
$newDetail.Range("H1").Value = "仓位排名"

$newDetail.Range("A2").Value = 0

# Columns B, D-G of the data row are numeric-looking but stored as TEXT
# (not numbers) in the source workbook - same convention as "2021-Q4".
# Force text storage via NumberFormat "@" while assigning, then immediately
# restore the plain default formatting (copied from a blank cell) so the
# cells end up as plain text values with no stray "@"/quote-prefix style,
# matching "2021-Q4"'s B2:G2 (t="s" / t="inlineStr", no explicit "s" attr).
$newDetail.Range("B2").NumberFormat = "@"
$newDetail.Range("D2:G2").NumberFormat = "@"
$newDetail.Range("B2").Value = "968013"
$newDetail.Range("C2").Value = "施罗德亚洲高息股债基金M"
$newDetail.Range("D2").Value = "297.64"
$newDetail.Range("E2").Value = "57.54"
$newDetail.Range("F2").Value = "1.26"
$newDetail.Range("G2").Value = "3.7503"

$blank.Copy()
$newDetail.Range("B2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newDetail.Range("H2").Value = 8

# --- Add the new "总计" sheet straight after "2022-Q1" ---------------------

$newTotals = $wb.Worksheets.Add($null, $newDetail)
$newTotals.Name = "总计"

# Seed formatting for the (3-row) totals table by copying the still-intact
# "2021-Q4" sheetPr/header/index-column styling twice: once for row 2, once
# for row 3 (the source table itself only has 2 rows).
$detail2021.Range("A1:D2").Copy()
$newTotals.Range("A1:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newTotals.Range("A2:D2").Copy()
$newTotals.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newTotals.Range("B1").Value = $totalsHeaderB
$newTotals.Range("C1").Value = $totalsHeaderC
$newTotals.Range("D1").Value = $totalsHeaderD

$newTotals.Range("A2").Value = 0
$newTotals.Range("B2").Value = "2022-Q1"
$newTotals.Range("C2").Value = 1
$newTotals.Range("D2").Value = 3.75

$newTotals.Range("A3").Value = 1
$newTotals.Range("B3").Value = $oldTotalsDate
$newTotals.Range("C3").Value = $oldTotalsCount
$newTotals.Range("D3").Value = $oldTotalsValue

$wb.Worksheets.Item("2021-Q4").Activate()
